$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RailInventory")

# Row 1 headers (M1:BL1)
$ws.Range("M1").Value = "n4OrphanContainer"
$ws.Range("N1").Value = "n4OrphanIsocode"
$ws.Range("O1").Value = "n4OrphanSlot"
$ws.Range("P1").Value = "n4OrphanUnit1"
$ws.Range("Q1").Value = "n4OrphanUnit2"
$ws.Range("R1").Value = "n4Seal"
$ws.Range("S1").Value = "n4YardLoc"
$ws.Range("T1").Value = "n4Lineoperator"
$ws.Range("U1").Value = "n4railFreightKind"
$ws.Range("V1").Value = "n4TankRails"
$ws.Range("W1").Value = "n4RailDirection"
$ws.Range("X1").Value = "n4RailNotes"
$ws.Range("Y1").Value = "n4DamageComp"
$ws.Range("Z1").Value = "n4DamageType"
$ws.Range("AA1").Value = "n4DamageSeverity"
$ws.Range("AB1").Value = "n4DamageLength"
$ws.Range("AC1").Value = "n4DamageWidth"
$ws.Range("AD1").Value = "n4DamageQuatity"
$ws.Range("AE1").Value = "n4DamageDeep"
$ws.Range("AF1").Value = "n4DamageLocation"
$ws.Range("AG1").Value = "n4RailInventoryPowerON"
$ws.Range("AH1").Value = "n4RailInventoryTemp"
$ws.Range("AI1").Value = "n4Gensetid"
$ws.Range("AJ1").Value = "n4railtype"
$ws.Range("AK1").Value = "n4OOGHeigth"
$ws.Range("AL1").Value = "n4OOGLeft"
$ws.Range("AM1").Value = "n4OOGRigth"
$ws.Range("AN1").Value = "n4OOGFlore"
$ws.Range("AO1").Value = "n4OOGUnits"
$ws.Range("AP1").Value = "n4RailInventoryPlacard"
$ws.Range("AQ1").Value = "n4Placard1"
$ws.Range("AR1").Value = "n4Placard2"
$ws.Range("AS1").Value = "n4Placard3"
$ws.Range("AT1").Value = "n4Placard4"
$ws.Range("AU1").Value = "n4Placard5"
$ws.Range("AV1").Value = "n4Placard6"
$ws.Range("AW1").Value = "n4RailInventoryGrd"
$ws.Range("AX1").Value = "n4RailInventorymaterial"
$ws.Range("AY1").Value = "n4RailInventoryS1"
$ws.Range("AZ1").Value = "n4RailInventoryS2"
$ws.Range("BA1").Value = "n4RailInventoryS3"
$ws.Range("BB1").Value = "n4RailInventoryS4"
$ws.Range("BC1").Value = "n4RailInventoryBundle"
$ws.Range("BD1").Value = "n4RailInventoryWeight"
$ws.Range("BE1").Value = "n4RailInventoryunit"
$ws.Range("BF1").Value = "n4RailInventoryTareWeight"
$ws.Range("BG1").Value = "n4RailInventorycscDate"
$ws.Range("BH1").Value = "n4RailInventoryMnf"
$ws.Range("BI1").Value = "n4RailInventorymrstatus"
$ws.Range("BJ1").Value = "n4RailInventoryGrdFood"
$ws.Range("BK1").Value = "n4RailInventorymaterialMade"
$ws.Range("BL1").Value = " "

# Row 2 data - text values
$ws.Range("M2").Value = "DINU17296000"
$ws.Range("P2").Value = "DINU1234567"
$ws.Range("Q2").Value = "DINU1231258"
$ws.Range("R2").Value = "s1"
$ws.Range("S2").Value = "1A"
$ws.Range("T2").Value = "ASW"
$ws.Range("U2").Value = "FCL"
$ws.Range("V2").Value = "Top"
$ws.Range("W2").Value = "Forward"
$ws.Range("Y2").Value = "DOOR"
$ws.Range("Z2").Value = "DAT"
$ws.Range("AA2").Value = "Major"
$ws.Range("AG2").Value = "NO"
$ws.Range("AH2").Value = "F"
$ws.Range("AP2").Value = "Yes"
$ws.Range("AQ2").Value = "CLASS 1.4"
$ws.Range("AR2").Value = "CLASS 1.4"
$ws.Range("AS2").Value = "CLASS 1.4"
$ws.Range("AT2").Value = "CLASS 1.4"
$ws.Range("AU2").Value = "CLASS 1.4"
$ws.Range("AV2").Value = "CLASS 1.4"
$ws.Range("AW2").Value = "FOOD"
$ws.Range("AX2").Value = "Steel"
$ws.Range("AY2").Value = "S12"
$ws.Range("AZ2").Value = "S2"
$ws.Range("BA2").Value = "S3"
$ws.Range("BB2").Value = "S4"
$ws.Range("BC2").Value = "ID1"
$ws.Range("BI2").Value = "MnRStatus"
$ws.Range("BJ2").Value = "FOOD"
$ws.Range("BK2").Value = "Steel"

# Row 2 data - numeric values
$ws.Range("N2").Value = 2200
$ws.Range("O2").Value = 2
$ws.Range("AB2").Value = 44
$ws.Range("AC2").Value = 11
$ws.Range("AD2").Value = 2
$ws.Range("AE2").Value = 15
$ws.Range("AJ2").Value = 9000
$ws.Range("AK2").Value = 11
$ws.Range("AL2").Value = 5
$ws.Range("AM2").Value = 10
$ws.Range("BD2").Value = 100
$ws.Range("BF2").Value = 2200
$ws.Range("BG2").Value = 12566
